$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.58598166666667
$ws.Range("H2").Value = 82.757945
$ws.Range("I2").Value = 0.2704460545904799
$ws.Range("J2").Value = 0.2704460545904799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.020559
$ws.Range("N2").Value = 90.061677
$ws.Range("O2").Value = 0.8829766276144534
$ws.Range("P2").Value = 0.8829766276144534
$ws.Range("Q2").Value = 828.1465901970852
$ws.Range("R2").Value = 7453.319311773766
$ws.Range("S2").Value = 0.2387975452339363
$ws.Range("T2").Value = 0.2387975452339363

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.58598166666667
$ws.Range("H3").Value = 82.757945
$ws.Range("I3").Value = 0.2704460545904799
$ws.Range("J3").Value = 0.2704460545904799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.426936666666667
$ws.Range("N3").Value = 4.28081
$ws.Range("O3").Value = 0.04196962907162197
$ws.Range("P3").Value = 0.04196962907162197
$ws.Range("Q3").Value = 39.36344872616111
$ws.Range("R3").Value = 354.27103853545
$ws.Range("S3").Value = 0.01135052059504607
$ws.Range("T3").Value = 0.01135052059504607

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.58598166666667
$ws.Range("H4").Value = 82.757945
$ws.Range("I4").Value = 0.2704460545904799
$ws.Range("J4").Value = 0.2704460545904799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.067853
$ws.Range("N4").Value = 3.203559
$ws.Range("O4").Value = 0.03140811737476231
$ws.Range("P4").Value = 0.0314081173747623
$ws.Range("Q4").Value = 29.45777328069501
$ws.Range("R4").Value = 265.1199595262551
$ws.Range("S4").Value = 0.00849420142611917
$ws.Range("T4").Value = 0.008494201426119166

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 27.58598166666667
$ws.Range("H5").Value = 82.757945
$ws.Range("I5").Value = 0.2704460545904799
$ws.Range("J5").Value = 0.2704460545904799
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.483919333333333
$ws.Range("N5").Value = 4.451758
$ws.Range("O5").Value = 0.04364562593916237
$ws.Range("P5").Value = 0.04364562593916237
$ws.Range("Q5").Value = 40.93537152414556
$ws.Range("R5").Value = 368.41834371731
$ws.Range("S5").Value = 0.01180378733537837
$ws.Range("T5").Value = 0.01180378733537837

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 26.23504533333333
$ws.Range("H6").Value = 78.705136
$ws.Range("I6").Value = 0.2572018131577233
$ws.Range("J6").Value = 0.2572018131577233
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.020559
$ws.Range("N6").Value = 90.061677
$ws.Range("O6").Value = 0.8829766276144534
$ws.Range("P6").Value = 0.8829766276144534
$ws.Range("Q6").Value = 787.590726297008
$ws.Range("R6").Value = 7088.316536673072
$ws.Range("S6").Value = 0.2271031895983293
$ws.Range("T6").Value = 0.2271031895983293

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 26.23504533333333
$ws.Range("H7").Value = 78.705136
$ws.Range("I7").Value = 0.2572018131577233
$ws.Range("J7").Value = 0.2572018131577233
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.426936666666667
$ws.Range("N7").Value = 4.28081
$ws.Range("O7").Value = 0.04196962907162197
$ws.Range("P7").Value = 0.04196962907162197
$ws.Range("Q7").Value = 37.43574813779555
$ws.Range("R7").Value = 336.9217332401599
$ws.Range("S7").Value = 0.01079466469477827
$ws.Range("T7").Value = 0.01079466469477827

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 26.23504533333333
$ws.Range("H8").Value = 78.705136
$ws.Range("I8").Value = 0.2572018131577233
$ws.Range("J8").Value = 0.2572018131577233
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.067853
$ws.Range("N8").Value = 3.203559
$ws.Range("O8").Value = 0.03140811737476231
$ws.Range("P8").Value = 0.0314081173747623
$ws.Range("Q8").Value = 28.015171864336
$ws.Range("R8").Value = 252.136546779024
$ws.Range("S8").Value = 0.00807822473665946
$ws.Range("T8").Value = 0.008078224736659456

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 26.23504533333333
$ws.Range("H9").Value = 78.705136
$ws.Range("I9").Value = 0.2572018131577233
$ws.Range("J9").Value = 0.2572018131577233
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.483919333333333
$ws.Range("N9").Value = 4.451758
$ws.Range("O9").Value = 0.04364562593916237
$ws.Range("P9").Value = 0.04364562593916237
$ws.Range("Q9").Value = 38.93069098100977
$ws.Range("R9").Value = 350.376218829088
$ws.Range("S9").Value = 0.01122573412795632
$ws.Range("T9").Value = 0.01122573412795632

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.31506333333334
$ws.Range("H10").Value = 87.94519000000001
$ws.Range("I10").Value = 0.2873975381543141
$ws.Range("J10").Value = 0.2873975381543141
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.020559
$ws.Range("N10").Value = 90.061677
$ws.Range("O10").Value = 0.8829766276144534
$ws.Range("P10").Value = 0.8829766276144534
$ws.Range("Q10").Value = 880.0545883870702
$ws.Range("R10").Value = 7920.491295483631
$ws.Range("S10").Value = 0.2537653090241925
$ws.Range("T10").Value = 0.2537653090241925

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 29.31506333333334
$ws.Range("H11").Value = 87.94519000000001
$ws.Range("I11").Value = 0.2873975381543141
$ws.Range("J11").Value = 0.2873975381543141
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.426936666666667
$ws.Range("N11").Value = 4.28081
$ws.Range("O11").Value = 0.04196962907162197
$ws.Range("P11").Value = 0.04196962907162197
$ws.Range("Q11").Value = 41.83073875598889
$ws.Range("R11").Value = 376.4766488039
$ws.Range("S11").Value = 0.01206196807243389
$ws.Range("T11").Value = 0.01206196807243389

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 29.31506333333334
$ws.Range("H12").Value = 87.94519000000001
$ws.Range("I12").Value = 0.2873975381543141
$ws.Range("J12").Value = 0.2873975381543141
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.067853
$ws.Range("N12").Value = 3.203559
$ws.Range("O12").Value = 0.03140811737476231
$ws.Range("P12").Value = 0.0314081173747623
$ws.Range("Q12").Value = 31.30417832569001
$ws.Range("R12").Value = 281.7376049312101
$ws.Range("S12").Value = 0.009026615611568428
$ws.Range("T12").Value = 0.009026615611568426

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 29.31506333333334
$ws.Range("H13").Value = 87.94519000000001
$ws.Range("I13").Value = 0.2873975381543141
$ws.Range("J13").Value = 0.2873975381543141
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.483919333333333
$ws.Range("N13").Value = 4.451758
$ws.Range("O13").Value = 0.04364562593916237
$ws.Range("P13").Value = 0.04364562593916237
$ws.Range("Q13").Value = 43.50118923822446
$ws.Range("R13").Value = 391.5107031440201
$ws.Range("S13").Value = 0.01254364544611934
$ws.Range("T13").Value = 0.01254364544611934

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 18.86569966666666
$ws.Range("H14").Value = 56.59709899999999
$ws.Range("I14").Value = 0.1849545940974826
$ws.Range("J14").Value = 0.1849545940974826
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.020559
$ws.Range("N14").Value = 90.061677
$ws.Range("O14").Value = 0.8829766276144534
$ws.Range("P14").Value = 0.8829766276144534
$ws.Range("Q14").Value = 566.3588499194469
$ws.Range("R14").Value = 5097.229649275023
$ws.Range("S14").Value = 0.1633105837579953
$ws.Range("T14").Value = 0.1633105837579953

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 18.86569966666666
$ws.Range("H15").Value = 56.59709899999999
$ws.Range("I15").Value = 0.1849545940974826
$ws.Range("J15").Value = 0.1849545940974826
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.426936666666667
$ws.Range("N15").Value = 4.28081
$ws.Range("O15").Value = 0.04196962907162197
$ws.Range("P15").Value = 0.04196962907162197
$ws.Range("Q15").Value = 26.92015859668777
$ws.Range("R15").Value = 242.28142737019
$ws.Range("S15").Value = 0.007762475709363748
$ws.Range("T15").Value = 0.007762475709363749

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 18.86569966666666
$ws.Range("H16").Value = 56.59709899999999
$ws.Range("I16").Value = 0.1849545940974826
$ws.Range("J16").Value = 0.1849545940974826
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.067853
$ws.Range("N16").Value = 3.203559
$ws.Range("O16").Value = 0.03140811737476231
$ws.Range("P16").Value = 0.0314081173747623
$ws.Range("Q16").Value = 20.145793986149
$ws.Range("R16").Value = 181.312145875341
$ws.Range("S16").Value = 0.005809075600415255
$ws.Range("T16").Value = 0.005809075600415253

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 18.86569966666666
$ws.Range("H17").Value = 56.59709899999999
$ws.Range("I17").Value = 0.1849545940974826
$ws.Range("J17").Value = 0.1849545940974826
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.483919333333333
$ws.Range("N17").Value = 4.451758
$ws.Range("O17").Value = 0.04364562593916237
$ws.Range("P17").Value = 0.04364562593916237
$ws.Range("Q17").Value = 27.99517647222689
$ws.Range("R17").Value = 251.956588250042
$ws.Range("S17").Value = 0.008072459029708335
$ws.Range("T17").Value = 0.008072459029708335
